# "This is my first commit" -- add a new learner (Parul Chaubey) to Sheet1,
# mirroring the existing rows (Name / EmailID with a mailto: hyperlink),
# then leave the workbook focused back on Sheet1 at the newly-added row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Add the new row to Sheet1 (Name / EmailID / Skills / Team table) ---
$ws1.Range("A6").Value = "Parul Chaubey"
$ws1.Range("B6").Value = "parulchaubey88@gmail.com"

# Hyperlink the new email address, matching the mailto: links used for
# every other row in this column.
$ws1.Hyperlinks.Add($ws1.Range("B6"), "mailto:parulchaubey88@gmail.com")

# Give the new email cell the same "Hyperlink" look as the rest of column B.
$ws1.Range("B6").Style = "Hyperlink"

# --- Reset Sheet2's selection back to its default (A1) ---
[void]$ws2.Range("A1").Select()

# --- Make Sheet1 the active sheet again, selecting the cell right after
#     the newly-typed entry (mirrors where the author's cursor ended up). ---
[void]$ws1.Activate()
[void]$ws1.Range("C6").Select()
